$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    3 = @(0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.719096746035642)
    4 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    5 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    6 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.21235251628842)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    8 = @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 1, 4.075514443323626)
    9 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.034748368925986)
    10 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    11 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698)
    12 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    13 = @(0.3464964993005633, 0.3375848360084654, 16.98373111632243, 0.4998867070740569, 0, 18.16769915870551)
    14 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    15 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
    16 = @(0.7287194209349384, 0.05231270169004087, 0.1529057820181812, 6.48142807727062, 1, 7.415365981913781)
    17 = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 2.27892381503245)
    18 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    19 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 1, 4.266567469086664)
    20 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698)
    21 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    22 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245)
    23 = @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 0, 4.075514443323626)
    24 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    25 = @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 2.652525447291612)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # B=2 ... G=7
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
